$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose textual content looks like a number or percentage need an
# explicit Text number format first, otherwise Excel auto-converts the
# assigned string into a numeric/percentage value instead of keeping it as text.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "309.65"
Set-TextValue $ws.Range("E2") "-0.98%"
Set-TextValue $ws.Range("D3") "37.11"
Set-TextValue $ws.Range("E3") "-2.15%"
Set-TextValue $ws.Range("D4") "5.125"
Set-TextValue $ws.Range("E4") "-0.27%"
Set-TextValue $ws.Range("D5") "0.07798"
Set-TextValue $ws.Range("E5") "-1.56%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D6") "4.401"
Set-TextValue $ws.Range("E6") "-0.11%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D7") "8.309"
Set-TextValue $ws.Range("E7") "0.73%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D8") "1.863"
Set-TextValue $ws.Range("E8") "-3.15%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.949"
Set-TextValue $ws.Range("E9") "4.56%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D10") "0.9249"
Set-TextValue $ws.Range("E10") "-0.25%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.1138"
Set-TextValue $ws.Range("E11") "-5.58%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D12") "0.1883"
Set-TextValue $ws.Range("E12") "-2.35%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D13") "0.08899"
Set-TextValue $ws.Range("E13") "-3.97%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03300"
Set-TextValue $ws.Range("E14") "-1.43%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09580"
Set-TextValue $ws.Range("E15") "-0.63%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001381"
Set-TextValue $ws.Range("E16") "0.76%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.006206"
Set-TextValue $ws.Range("E17") "4.79%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.394"
Set-TextValue $ws.Range("E18") "-4.17%"
Set-TextValue $ws.Range("D19") "0.3450"
Set-TextValue $ws.Range("E19") "0.11%"
Set-TextValue $ws.Range("D20") "6.399"
Set-TextValue $ws.Range("E20") "20.94%"
Set-TextValue $ws.Range("D21") "0.1292"
Set-TextValue $ws.Range("E21") "0.58%"
Set-TextValue $ws.Range("D22") "0.2374"
Set-TextValue $ws.Range("E22") "-8.30%"
Set-TextValue $ws.Range("D23") "0.04343"
Set-TextValue $ws.Range("E23") "-0.78%"
Set-TextValue $ws.Range("E24") "-4.03%"
Set-TextValue $ws.Range("D25") "0.004271"
Set-TextValue $ws.Range("E25") "-0.36%"
Set-TextValue $ws.Range("D26") "0.0001402"
Set-TextValue $ws.Range("E26") "7.89%"
Set-TextValue $ws.Range("D27") "0.0002905"
Set-TextValue $ws.Range("D39") "0.02130"
Set-TextValue $ws.Range("E39") "0.72%"
Set-TextValue $ws.Range("D40") "0.04984"
Set-TextValue $ws.Range("E40") "-2.34%"
Set-TextValue $ws.Range("D41") "0.007583"
Set-TextValue $ws.Range("E41") "-0.40%"
Set-TextValue $ws.Range("D42") "0.1354"
Set-TextValue $ws.Range("E42") "-0.32%"
Set-TextValue $ws.Range("D43") "0.008519"
Set-TextValue $ws.Range("E43") "-6.65%"
Set-TextValue $ws.Range("D44") "0.002074"
Set-TextValue $ws.Range("E44") "1.18%"
Set-TextValue $ws.Range("D45") "0.007986"
Set-TextValue $ws.Range("E45") "-7.30%"
Set-TextValue $ws.Range("D46") "0.00006582"
Set-TextValue $ws.Range("E46") "-1.61%"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "0.19%"
Set-TextValue $ws.Range("E48") "13.76%"
Set-TextValue $ws.Range("E49") "20.49%"
Set-TextValue $ws.Range("E50") "0.19%"
Set-TextValue $ws.Range("E51") "0.19%"
